# Repull data, push all data, mean calculation
# Update the "dSF" column (F) values for the affected rows to reflect
# re-pulled / re-pushed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    5  = 2
    9  = 4
    11 = -4
    14 = 4
    15 = -3
    19 = -1
    20 = -1
    22 = -1
    23 = 4
    27 = 0
    29 = 2
    34 = 1
    38 = 4
    44 = -8
    46 = -3
    47 = -1
    53 = -4
    55 = 0
    57 = 0
    58 = -4
    59 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
